# Dataset spec sheet "etapaitem": add two new leading columns
# (convocatoria_ocds_id, item_etapa_id) ahead of the existing
# item_id / item_descripcion / item_clasificacion / item_cantidad / item_unidad
# headers, and replace the placeholder second row with real sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A:B columns (and everything after) two columns to the
# right so the current headers/data land on C:G, keeping their formatting.
$ws.Range("A:B").Insert(-4161) | Out-Null

# New header cells in the freed-up A1:B1 slot, formatted like the rest of
# the header row (bold font, thin border, centered/top alignment).
$ws.Range("A1").Value = "convocatoria_ocds_id"
$ws.Range("B1").Value = "item_etapa_id"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 2 sample data. B2 is intentionally left blank (no value for
# item_etapa_id in this sample row).
$ws.Range("A2").Value = "ocds-twb234-0005"
$ws.Range("D2").Value = "Servicio de consultoria"
$ws.Range("E2").Value = "Soporte y mantenimiento de hardware"
$ws.Range("G2").Value = "Cantidad"

# C2 / F2 hold numeric-looking values ("3245" / "1") that must stay text,
# matching the original sheet's string-typed cells. A leading apostrophe
# forces text entry; resetting the style back to Normal afterwards drops
# the quote-prefix formatting so the cells keep the sheet's default look.
$ws.Range("C2").Value = "'3245"
$ws.Range("C2").Style = "Normal"
$ws.Range("F2").Value = "'1"
$ws.Range("F2").Style = "Normal"

# Page margins match the regenerated spec export (inches -> points).
$ws.PageSetup.LeftMargin = 50.4
$ws.PageSetup.RightMargin = 50.4
$ws.PageSetup.TopMargin = 54.0
$ws.PageSetup.BottomMargin = 54.0
$ws.PageSetup.HeaderMargin = 21.599999999999998
$ws.PageSetup.FooterMargin = 21.599999999999998
